# Add a reprint row for "FOCUS sax mid" to the "July 2018" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

$row = 19

$ws.Cells.Item($row, 1).Value = "26-07-2018"
$ws.Cells.Item($row, 2).Value = "26-07-2018"
$ws.Cells.Item($row, 3).Value = "FOCUS sax 100 Scale mid"
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = "Polylite"
$ws.Cells.Item($row, 6).Value = 2
$ws.Cells.Item($row, 7).Value = 20
$ws.Cells.Item($row, 8).Value = 0.2
$ws.Cells.Item($row, 9).Value = "NA"

$ws.Range("A" + $row + ":I" + $row).HorizontalAlignment = -4108

$ws.Range("B19").Select()
